# Applies the "added 4wk low sales check" commit:
# - Forecast Comparison sheet: updated MyForecast (D), Trend (G),
#   Inventory Coverage (H), Stockout Risk (I), Reorder Urgency (J) and
#   Seasonality Index (L) for the weekly forecast rows.
# - Summary sheet: updated the aggregate forecast totals/extremes that
#   depend on the recalculated weekly forecast.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---------------------------------------

# Row 2 (W10)
$ws1.Cells.Item(2, "D").Value = 437
$ws1.Cells.Item(2, "G").Value = "High Volume Season"
$ws1.Cells.Item(2, "H").Value = 3.9
$ws1.Cells.Item(2, "L").Value = 1.16

# Row 3 (W11)
$ws1.Cells.Item(3, "D").Value = 1049
$ws1.Cells.Item(3, "G").Value = "High Volume Season"
$ws1.Cells.Item(3, "H").Value = 1.21
$ws1.Cells.Item(3, "L").Value = 0.99

# Row 4 (W12)
$ws1.Cells.Item(4, "D").Value = 1611
$ws1.Cells.Item(4, "G").Value = "High Volume Season"
$ws1.Cells.Item(4, "H").Value = 0.14
$ws1.Cells.Item(4, "I").Value = "High"
$ws1.Cells.Item(4, "J").Value = "Urgent"
$ws1.Cells.Item(4, "L").Value = 0.9

# Row 5 (W13)
$ws1.Cells.Item(5, "D").Value = 1413
$ws1.Cells.Item(5, "G").Value = "High Volume Season"
$ws1.Cells.Item(5, "H").Value = 0
$ws1.Cells.Item(5, "I").Value = "High"
$ws1.Cells.Item(5, "J").Value = "Urgent"
$ws1.Cells.Item(5, "L").Value = 0.86

# Row 6 (W14)
$ws1.Cells.Item(6, "D").Value = 357
$ws1.Cells.Item(6, "G").Value = "High Volume Season"
$ws1.Cells.Item(6, "H").Value = 0
$ws1.Cells.Item(6, "I").Value = "High"
$ws1.Cells.Item(6, "J").Value = "Urgent"
$ws1.Cells.Item(6, "L").Value = 0.95

# Row 7 (W15)
$ws1.Cells.Item(7, "D").Value = 164
$ws1.Cells.Item(7, "G").Value = "High Volume Season"
$ws1.Cells.Item(7, "H").Value = 0
$ws1.Cells.Item(7, "I").Value = "High"
$ws1.Cells.Item(7, "J").Value = "Urgent"
$ws1.Cells.Item(7, "L").Value = 0.86

# Row 8 (W16)
$ws1.Cells.Item(8, "D").Value = 164
$ws1.Cells.Item(8, "G").Value = "High Volume Season"
$ws1.Cells.Item(8, "H").Value = 0
$ws1.Cells.Item(8, "I").Value = "High"
$ws1.Cells.Item(8, "L").Value = 1.06

# Row 9 (W17)
$ws1.Cells.Item(9, "D").Value = 805
$ws1.Cells.Item(9, "G").Value = "High Volume Season"
$ws1.Cells.Item(9, "L").Value = 1.05

# Row 10 (W18)
$ws1.Cells.Item(10, "D").Value = 1733
$ws1.Cells.Item(10, "G").Value = "High Volume Season"
$ws1.Cells.Item(10, "L").Value = 1.14

# Row 11 (W19)
$ws1.Cells.Item(11, "D").Value = 1115
$ws1.Cells.Item(11, "G").Value = "High Volume Season"
$ws1.Cells.Item(11, "L").Value = 1.16

# Row 12 (W20)
$ws1.Cells.Item(12, "D").Value = 164
$ws1.Cells.Item(12, "G").Value = "High Volume Season"
$ws1.Cells.Item(12, "L").Value = 1.04

# Row 13 (W21)
$ws1.Cells.Item(13, "D").Value = 164
$ws1.Cells.Item(13, "G").Value = "High Volume Season"
$ws1.Cells.Item(13, "L").Value = 0.93

# Row 14 (W22)
$ws1.Cells.Item(14, "D").Value = 164
$ws1.Cells.Item(14, "G").Value = "High Volume Season"
$ws1.Cells.Item(14, "L").Value = 1.11

# Row 15 (W23)
$ws1.Cells.Item(15, "D").Value = 164
$ws1.Cells.Item(15, "G").Value = "High Volume Season"
$ws1.Cells.Item(15, "L").Value = 1.18

# Row 16 (W24)
$ws1.Cells.Item(16, "D").Value = 164
$ws1.Cells.Item(16, "G").Value = "High Volume Season"
$ws1.Cells.Item(16, "L").Value = 1.11

# Row 17 (W25)
$ws1.Cells.Item(17, "D").Value = 164
$ws1.Cells.Item(17, "G").Value = "High Volume Season"
$ws1.Cells.Item(17, "L").Value = 0.82

# --- Summary sheet -----------------------------------------------------
# These cells are stored as text, so force text formatting before writing
# the numeric-looking strings to keep them as text values.

$ws2.Range("B9").NumberFormat  = "@"
$ws2.Range("B9").Value  = "9845"   # Total Forecast (16 Weeks)

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "6005"   # Total Forecast (8 Weeks)

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "4512"   # Total Forecast (4 Weeks)

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "1734"   # Max Forecast

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "165"    # Min Forecast
